# Generate Report for Archive
#
# 1) Shared string used for the localization status cells changes from
#    "Ready for handoff" to "In Translation" (it's the same shared-string
#    entry referenced by every occurrence, so re-writing the cell values
#    updates every usage: the "zh-cn"/"de-de" columns on the Overview sheet,
#    and the "Status" column on each language sheet).
# 2) The "Status"-ish columns get narrower (their column width shrinks).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (E) / de-de (F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
